# Apply "Renamed few transcripts. Updated the DataSheet" edits.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell F1: "Teacher Tag" -> "T Tag"
$ws.Range("F1").Value = "T Tag"

# Rows where column D currently holds "T/R1" -> rename to "T"
$rowsTR1 = @(2,3,4,29,30,33,39,44,46,48,49,50,51,52,54,57,60,62,63,64,66,67)
foreach ($r in $rowsTR1) {
    $ws.Cells.Item($r, 4).Value = "T"
}

# Rows where column D currently holds "T/R2" -> rename to "T 2"
$rowsTR2 = @(70,72,75,76,77,79,81,85,87,88,89,91,93,95,97,99,101,102,104,107,109,110,112,113)
foreach ($r in $rowsTR2) {
    $ws.Cells.Item($r, 4).Value = "T 2"
}
